$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, pushing existing rows 9-14 down to 10-15.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the Scientific Reports entry (column E only).
$ws.Range("E9").Value = "Scientific Reports"

# Adjust column D width (it no longer needs to fit the long theme-issue text
# that used to live in row 13/now row 14 once other columns reflow).
$ws.Columns.Item(4).ColumnWidth = 21.5703125

# Restore the active selection to E9, matching the saved view state.
$ws.Range("E9").Select()
